$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C5").Value = "wait_seconds"
$ws.Range("D5").Value = 2
$ws.Range("D11").Value = 3
$ws.Range("D5").Select()
